# Insert a new weekly price record for Vega Monumental Concepción - Betarraga.
# This pushes the existing rows 265..321 down to 266..322 and populates the
# newly opened row 265 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 265 (existing row 265 and everything below it
# shift down by one row).
$ws.Rows.Item(265).Insert()

# Fill in the new row 265 with the new weekly record.
$ws.Cells.Item(265, 1).Value = 11
$ws.Cells.Item(265, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(265, 3).Value = "Bíobío"
$ws.Cells.Item(265, 4).Value = 44722
$ws.Cells.Item(265, 5).Value = 8
$ws.Cells.Item(265, 6).Value = 100114014
$ws.Cells.Item(265, 7).Value = "Betarraga"
$ws.Cells.Item(265, 8).Value = "Sin especificar"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 1150
$ws.Cells.Item(265, 11).Value = 600
$ws.Cells.Item(265, 12).Value = 650
$ws.Cells.Item(265, 13).Value = 624
$ws.Cells.Item(265, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(265, 15).Value = "Región Metropolitana"
$ws.Cells.Item(265, 16).Value = 125
$ws.Cells.Item(265, 17).Value = 5
$ws.Cells.Item(265, 18).Value = "Hortaliza"
